# Extend the work-log table (Tabelle1) with three new log entries
# (rows 55-57), which were previously blank placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 55-57 were blank placeholder rows (style "1" only, no number
# format). Re-use the existing date/time formatting from the row above
# (row 54: date style, time style, time style, left-aligned text style)
# by copying its formats down, so the new cells pick up the same style
# indexes already present in the workbook instead of creating new ones.
$ws.Range("B54:E54").Copy()
$ws.Range("B55:E57").PasteSpecial(-4122)   # xlPasteFormats

# C56 holds free text ("09:30?") instead of a time value, and keeps the
# plain general/centered style (same as header-row style), not the time
# style, so fix that cell's format back to the plain style.
$ws.Range("B1").Copy()
$ws.Range("C56").PasteSpecial(-4122)       # xlPasteFormats

# --- Row 55: 20-11-2022 -------------------------------------------------
$ws.Range("B55").Value = 44885
$ws.Range("C55").Value = 0.375
$ws.Range("D55").Value = 0.77083333333333337
$ws.Range("E55").Value = "designed v0.2 of pcb"

# --- Row 56: 21-11-2022 -------------------------------------------------
$ws.Range("B56").Value = 44886
$ws.Range("C56").Value = "09:30?"
$ws.Range("D56").Value = 0.70833333333333337
$ws.Range("E56").Value = "finished designing v0.2 of pcb"

# --- Row 57: 22-11-2022 -------------------------------------------------
$ws.Range("B57").Value = 44887
$ws.Range("C57").Value = 0.375
$ws.Range("D57").Value = 0.70833333333333337
$ws.Range("E57").Value = "Attended UDC lessons and gave students feedback. Created pcb heater design and worked on camera. "

# Move the selection/active cell to the last entry and scroll the sheet
# so that the new rows are visible, matching where the user ended up.
[void]$ws.Range("E57").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 2
